$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-12 21:48:51"
$ws.Range("E3").Value = "2026-02-12 21:48:54"
$ws.Range("E4").Value = "2026-02-12 21:48:57"
$ws.Range("J4").Value = "999.6 hPa"
$ws.Range("E5").Value = "2026-02-12 21:49:00"
$ws.Range("E6").Value = "2026-02-12 21:49:02"
$ws.Range("J6").Value = "999.5 hPa"
$ws.Range("O6").Value = "15.8 °C"
$ws.Range("E7").Value = "2026-02-12 21:49:05"
$ws.Range("J7").Value = "1002.2 hPa"
$ws.Range("E8").Value = "2026-02-12 21:49:08"
$ws.Range("J8").Value = "1001.6 hPa"
$ws.Range("O8").Value = "13.0 °C"
$ws.Range("E9").Value = "2026-02-12 21:49:11"
$ws.Range("H9").Value = "'65%"
$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("O9").Value = "12.9 °C"
$ws.Range("E10").Value = "2026-02-12 21:49:14"
$ws.Range("H10").Value = "'47%"
$ws.Range("G10").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("E11").Value = "2026-02-12 21:49:16"
$ws.Range("O11").Value = "9.3 °C"
$ws.Range("E12").Value = "2026-02-12 21:49:19"
$ws.Range("O12").Value = "12.6 °C"
$ws.Range("E13").Value = "2026-02-12 21:49:21"
$ws.Range("J13").Value = "1002.1 hPa"
$ws.Range("O13").Value = "7.6 °C"
$ws.Range("E14").Value = "2026-02-12 21:49:24"
$ws.Range("H14").Value = "'36%"
$ws.Range("G14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("N14").Value = "13.7 °C 21:28 TU"
$ws.Range("E15").Value = "2026-02-12 21:49:26"
$ws.Range("H15").Value = "'52%"
$ws.Range("G15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("E16").Value = "2026-02-12 21:49:29"
$ws.Range("E17").Value = "2026-02-12 21:49:31"
$ws.Range("E18").Value = "2026-02-12 21:49:34"
$ws.Range("J18").Value = "999.9 hPa"
$ws.Range("N18").Value = "12.1 °C 21:29 TU"
$ws.Range("O18").Value = "16.7 °C"
$ws.Range("E19").Value = "2026-02-12 21:49:37"
$ws.Range("O19").Value = "8.0 °C"
$ws.Range("E20").Value = "2026-02-12 21:49:40"
$ws.Range("E21").Value = "2026-02-12 21:49:42"
$ws.Range("H21").Value = "'49%"
$ws.Range("G21").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("J21").Value = "1002.6 hPa"
$ws.Range("N21").Value = "4.7 °C 21:17 TU"
$ws.Range("O21").Value = "9.1 °C"
$ws.Range("E22").Value = "2026-02-12 21:49:45"
$ws.Range("O22").Value = "-5.7 °C"
$ws.Range("E23").Value = "2026-02-12 21:49:48"
$ws.Range("E24").Value = "2026-02-12 21:49:51"
$ws.Range("J24").Value = "1006.8 hPa"
$ws.Range("E25").Value = "2026-02-12 21:49:54"
$ws.Range("E26").Value = "2026-02-12 21:49:56"
$ws.Range("J26").Value = "999.2 hPa"
$ws.Range("E27").Value = "2026-02-12 21:49:59"
$ws.Range("O27").Value = "-1.6 °C"
$ws.Range("E28").Value = "2026-02-12 21:50:02"
$ws.Range("J28").Value = "999.3 hPa"
$ws.Range("K28").Value = "13.2 MJ/m2"
$ws.Range("O28").Value = "13.8 °C"
$ws.Range("E29").Value = "2026-02-12 21:50:05"
$ws.Range("H29").Value = "'60%"
$ws.Range("G29").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("N29").Value = "5.9 °C 21:27 TU"
$ws.Range("O29").Value = "14.2 °C"
$ws.Range("E30").Value = "2026-02-12 21:50:07"
$ws.Range("J30").Value = "999.7 hPa"
$ws.Range("N30").Value = "7.2 °C 21:29 TU"
$ws.Range("O30").Value = "12.0 °C"
$ws.Range("E31").Value = "2026-02-12 21:50:10"
$ws.Range("J31").Value = "999.2 hPa"
$ws.Range("O31").Value = "14.2 °C"
$ws.Range("E32").Value = "2026-02-12 21:50:12"
$ws.Range("O32").Value = "8.0 °C"
$ws.Range("E33").Value = "2026-02-12 21:50:15"
$ws.Range("J33").Value = "1001.8 hPa"
$ws.Range("E34").Value = "2026-02-12 21:50:18"
$ws.Range("E35").Value = "2026-02-12 21:50:21"
$ws.Range("E36").Value = "2026-02-12 21:50:23"
$ws.Range("J36").Value = "1000.0 hPa"
$ws.Range("O36").Value = "14.4 °C"
$ws.Range("E37").Value = "2026-02-12 21:50:26"
$ws.Range("J37").Value = "1000.7 hPa"
$ws.Range("N37").Value = "4.5 °C 21:24 TU"
$ws.Range("O37").Value = "9.8 °C"
$ws.Range("E38").Value = "2026-02-12 21:50:28"
$ws.Range("E39").Value = "2026-02-12 21:50:31"
$ws.Range("E40").Value = "2026-02-12 21:50:34"
$ws.Range("J40").Value = "1003.4 hPa"
$ws.Range("N40").Value = "4.7 °C 21:18 TU"
$ws.Range("O40").Value = "9.4 °C"
$ws.Range("E41").Value = "2026-02-12 21:50:37"
$ws.Range("J41").Value = "1005.8 hPa"
$ws.Range("O41").Value = "17.1 °C"
$ws.Range("E42").Value = "2026-02-12 21:50:39"
$ws.Range("N42").Value = "7.6 °C 21:27 TU"
$ws.Range("O42").Value = "13.9 °C"
$ws.Range("E43").Value = "2026-02-12 21:50:42"
$ws.Range("O43").Value = "11.9 °C"
$ws.Range("E44").Value = "2026-02-12 21:50:44"
$ws.Range("E45").Value = "2026-02-12 21:50:47"
$ws.Range("H45").Value = "'54%"
$ws.Range("G45").Copy()
$ws.Range("H45").PasteSpecial(-4122)
$ws.Range("J45").Value = "1005.3 hPa"
$ws.Range("E46").Value = "2026-02-12 21:50:50"
$ws.Range("J46").Value = "1007.5 hPa"
$ws.Range("O46").Value = "15.8 °C"
$excel.CutCopyMode = 0
